$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0,
    120.2000000000006,
    110,
    191.1999999999998,
    189.1999999999998,
    148.2000000000006,
    182.2000000000005,
    109,
    200,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
